# Update latest output (run 60)
# Applies updated optimisation results to the "Schedule" and "Detailed" sheets.

$wb = $excel.ActiveWorkbook
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsDetailed = $wb.Worksheets.Item("Detailed")

# --- Schedule sheet updates ---
$wsSchedule.Range("E4").Value = 216.116355
$wsSchedule.Range("F4").Value = 8.167662698412698
$wsSchedule.Range("A5").Value = 46040.3125
$wsSchedule.Range("B5").Value = 46040.8125
$wsSchedule.Range("E5").Value = -6.772457249999995
$wsSchedule.Range("F5").Value = -0.149304613095238

# --- Detailed sheet updates: Price (B) values ---
$wsDetailed.Range("B28").Value = 36.06045
$wsDetailed.Range("B29").Value = 31.3867
$wsDetailed.Range("B30").Value = 30.40942
$wsDetailed.Range("B31").Value = -9.99
$wsDetailed.Range("B32").Value = -6.8
$wsDetailed.Range("B33").Value = 22.07
$wsDetailed.Range("B34").Value = 36.06036
$wsDetailed.Range("B35").Value = 25.0354
$wsDetailed.Range("B36").Value = 4.53593
$wsDetailed.Range("B37").Value = -0.24301
$wsDetailed.Range("B38").Value = 0.00014
$wsDetailed.Range("B39").Value = 0.00014
$wsDetailed.Range("B40").Value = 17.07182
$wsDetailed.Range("B41").Value = 57.18142
$wsDetailed.Range("B43").Value = 57.18142
$wsDetailed.Range("B44").Value = 46.60383
$wsDetailed.Range("B45").Value = 46.91075
$wsDetailed.Range("B46").Value = 36.06045
$wsDetailed.Range("B47").Value = 47.4748
$wsDetailed.Range("B48").Value = 56.98
$wsDetailed.Range("B49").Value = 56.83
$wsDetailed.Range("B50").Value = 47.78559
$wsDetailed.Range("B52").Value = 57.06016
$wsDetailed.Range("B61").Value = 57.06003
$wsDetailed.Range("B64").Value = 36.0595
$wsDetailed.Range("B65").Value = 35.87994
$wsDetailed.Range("B66").Value = 0.63134
$wsDetailed.Range("B67").Value = 0.7
$wsDetailed.Range("B68").Value = 0.61554
$wsDetailed.Range("B69").Value = 0.59266
$wsDetailed.Range("B70").Value = 8.38785
$wsDetailed.Range("B71").Value = 27.10585
$wsDetailed.Range("B72").Value = 0.7
$wsDetailed.Range("B73").Value = 0.66306
$wsDetailed.Range("B74").Value = 0
$wsDetailed.Range("B75").Value = -5.21817
$wsDetailed.Range("B76").Value = -5.43091
$wsDetailed.Range("B77").Value = -5.57994
$wsDetailed.Range("B78").Value = -13.5
$wsDetailed.Range("B80").Value = -13.5
$wsDetailed.Range("B81").Value = -12.01
$wsDetailed.Range("B82").Value = -7.44788
$wsDetailed.Range("B83").Value = -8.045170000000001
$wsDetailed.Range("B84").Value = -7.66355
$wsDetailed.Range("B85").Value = -6.12644
$wsDetailed.Range("B86").Value = -5.85772
$wsDetailed.Range("B87").Value = 0.0001
$wsDetailed.Range("B88").Value = 22.15733
$wsDetailed.Range("B89").Value = 44.02345
$wsDetailed.Range("B90").Value = 47.60158
$wsDetailed.Range("B91").Value = 46.82896
$wsDetailed.Range("B92").Value = 46.12397
$wsDetailed.Range("B93").Value = 36.2
$wsDetailed.Range("B94").Value = 55.31303
$wsDetailed.Range("B95").Value = 36.0601

# --- Detailed sheet updates: Type (C) / Pump_Status (E) text values ---
$wsDetailed.Range("C31").Value = "historical"
$wsDetailed.Range("C32").Value = "historical"
$wsDetailed.Range("E64").Value = "OFF"
$wsDetailed.Range("E88").Value = "ON"

